# Weekly driver report update for 2025-04-28
# Updates the "Bad Drivers" and "Good Drivers" tables on the Driver Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152

# ---------------------------------------------------------------------------
# Section 1: "Bad Drivers" table (rows 3-8 before -> rows 3-6 after)
# Two rows worth of data were dropped (old 23.60.1.2 / 23.40.0.4 rows), and
# the remaining three driver rows + totals row were refreshed with new
# numbers for the week of 2025-04-28.
# ---------------------------------------------------------------------------

$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.3.2"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 549
$ws.Range("D3").Value = 98.8

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 610
$ws.Range("D4").Value = 98.8

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.240.0.6"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 1175
$ws.Range("D5").Value = 98.90000000000001

# Move the Totals row from row 8 up to row 6 with refreshed totals.
$ws.Range("A6").Value = "Totals:"
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 2334

# Clear out what used to be rows 7 and 8 (old totals row + the extra driver
# row) now that the table only spans rows 1-6.
$ws.Range("A7:J8").Clear()

# ---------------------------------------------------------------------------
# Section 2: "Good Drivers" table (rows 14-21 before -> rows 12-20 after)
# The whole block shifts up two rows (following the Bad Drivers deletions)
# and gains one new driver entry at the top of the list.
# ---------------------------------------------------------------------------

$ws.Range("A12").Value = "Good Drivers (Roaming > 99.8%)"

$ws.Range("A13").Value = "Adapter-Driver"
$ws.Range("B13").Value = "Total Samples"
$ws.Range("D13").Value = "Good Roaming Calculation (%)"
$ws.Range("E13").Value = "Driver Vintage"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B14").Value = 11128
$ws.Range("D14").Value = 100

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 486214
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2024-11-10"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 79953
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "'2021-08-18"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B17").Value = 35355
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2021-04-27"

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B18").Value = 65425
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2020-08-05"

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B19").Value = 117653
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2020-01-06"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B20").Value = 56018
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = "'2019-12-14"

# Apply number formats / alignment matching the rest of the column so the
# newly written cells look like their neighbours.
$ws.Range("B14:B20").NumberFormat = "#,##0"
$ws.Range("B14:B20").HorizontalAlignment = $xlRight
$ws.Range("D14:D20").HorizontalAlignment = $xlRight
$ws.Range("E14:E20").HorizontalAlignment = $xlRight

# Old row 21 (56018 / 21.60.2.1 entry) is now fully represented by row 20
# above, so clear the leftover row and the rest of the old tail rows.
$ws.Range("A21:J26").Clear()
